$wb = $excel.ActiveWorkbook

# --- Add the new "Storage CO2" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Storage CO2"

# --- Column width ---
$ws.Columns.Item(2).ColumnWidth = 39.5

# --- Content (written in the same order the shared-string table was built in) ---
$ws.Range("B6").Value = "offshore_ccs_potential_mt_per_year"
$ws.Range("C6").Value = 52
$ws.Range("C6").NumberFormat = "0"
$ws.Range("D6").Value = "Mton/year"

$ws.Range("B8").Value = "offshore_ccs_potential_mt_per_year"
$ws.Range("C8").Formula = "=C6*C7"
$ws.Range("C8").NumberFormat = "0"
$ws.Range("C8").Borders.Weight = -4138
$ws.Range("D8").Value = "Mton/year"

$ws.Range("B2").Value = "Based on the study ""Nationale CO2-opslagbehoefte tot 2035"" by Royal HaskoningDHV"

$ws.Range("B5").Value = "Yearly potential for offshore CO2 storage"

$ws.Range("B7").Value = "Correctection factor to allow new studies"
$ws.Range("C7").Value = 1.25

$ws.Range("C5").Value = 2035

$ws.Range("B3").Value = "Scenario: ""Maximale afvang"", year 2035"

# --- Selection / active state ---
$ws.Range("F9").Select()
$ws.Activate()
